# feat: prepare for accepting recipient without IE
#
# - Entidades (sheet1): drop the "[inscrição municipal]" placeholder/column
#   content, move "número" header from B to K, and insert new address
#   columns (cep, bairro, logradouro (tipo), logradouro (nome)) plus a new
#   "inscrição estadual" header in column B.
# - Dados das listas (sheet4): add a 9th column listing the logradouro
#   (tipo) options (Avenida / Rua / Estrada) and wire a new list
#   validation on Entidades!I2:I5 to it.
# - Make "Entidades" the active/selected tab again (was "Dados de Produtos
#   e Serviços NF").

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Entidades")
$ws4 = $wb.Worksheets.Item("Dados das listas")

# --- Entidades: drop the old "[inscrição municipal]" value -----------------
# Clearing this (and only this) reference lets the now-unused shared string
# get garbage collected on save, which is what re-packs the shared string
# table the same way the target workbook does.
$ws1.Range("B2").Value = ""

# --- Entidades: relocate "número" header from B1 to K1 ---------------------
$ws1.Range("K1").Value = "número"

# --- Entidades: new header cells (order matters for shared-string ids) -----
$ws1.Range("G1").Value = "cep"
$ws1.Range("H1").Value = "bairro"
$ws1.Range("I1").Value = "logradouro (tipo)"
$ws1.Range("J1").Value = "logradouro (nome)"

# --- Dados das listas: new "logradouro (tipo)" list (header + options) ----
$ws4.Range("I1").Value = "logradouro (tipo)"
$ws4.Range("I1").Font.Bold = $true
$ws4.Range("I2").Value = "Avenida"
$ws4.Range("I3").Value = "Rua"
$ws4.Range("I4").Value = "Estrada"

# --- Entidades: new "inscrição estadual" header replacing B1 "número" -----
$ws1.Range("B1").Value = "inscrição estadual"

# --- Column widths ----------------------------------------------------------
$ws1.Range("B1:C1").ColumnWidth = 16.333333333333332
$ws1.Range("E1:F1").ColumnWidth = 14.833333333333334
$ws1.Range("G1").ColumnWidth = 12.166666666666666
$ws1.Range("H1:J1").ColumnWidth = 19.5
$ws1.Range("K1").ColumnWidth = 8.333333333333334

$ws4.Range("I1").ColumnWidth = 15.451822916666666

# --- New data validation: Entidades!I2:I5 <- 'Dados das listas'!$I$2:$I$4 --
$dv = $ws1.Range("I2:I5").Validation
$dv.Add(3, 1, 1, "='Dados das listas'!`$I`$2:`$I`$4")
$dv.IgnoreBlank = $true
$dv.InCellDropdown = $true
$dv.ShowInput = $true
$dv.ShowError = $true

# --- Selections / active tab -----------------------------------------------
$ws4.Range("I1:I4").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("J6").Select() | Out-Null
